# Results updated against 1000 randomized hierarchies
#
# Updates the "Whole Genome" (row 5/6) and "Level 1" -> summary "No Overlap"
# rows (row 9/10) raw inputs on Sheet2 with the newly-computed values from
# the 1000-randomized-hierarchy run. The "Difference" rows (6 and 10) hold
# =B4-B5 style formulas, so they recalculate automatically once the new
# inputs are written.
#
# Also restores the view state (scroll position / active cell) that Calc
# persisted for each sheet after the edit.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2: row 5 ("R_M vs R_H->M" / Whole-Genome raw values) ---
$ws2.Range("B5").Value = 0.543135
$ws2.Range("C5").Value = 0.416325
$ws2.Range("D5").Value = 0.419712
$ws2.Range("E5").Value = 0.134032
$ws2.Range("F5").Value = 0.00922317
$ws2.Range("G5").Value = 0.0084504

# --- Sheet2: row 9 (No Overlap raw values) ---
$ws2.Range("B9").Value = 0.535098
$ws2.Range("C9").Value = 0.494541
$ws2.Range("D9").Value = 0.403448
$ws2.Range("E9").Value = 0.150666
$ws2.Range("F9").Value = 0.0140641
$ws2.Range("G9").Value = 0.0147694

# Rows 6 and 10 are "=<col>4-<col>5" / "=<col>8-<col>9" difference formulas,
# so they pick up the new values automatically on recalculation.

# --- Restore view/selection state ---
# Sheet1: scrolled so row 16 is at the top, with B46 as the active cell.
[void]$ws1.Range("A16").Select()
[void]$ws1.Range("B46").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1

# Sheet2: stays the active/tabbed sheet, active cell moves to G9.
[void]$ws2.Range("G9").Select()

$ws2.Activate()
